$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Update the time_taken timestamps in column F (rows 2-18) to reflect the
# re-run panel query time.
$data.Range("F2").Value = "2021-10-05 14:22:38.481370"
$data.Range("F3").Value = "2021-10-05 14:22:38.481378"
$data.Range("F4").Value = "2021-10-05 14:22:38.481381"
$data.Range("F5").Value = "2021-10-05 14:22:38.481383"
$data.Range("F6").Value = "2021-10-05 14:22:38.481386"
$data.Range("F7").Value = "2021-10-05 14:22:38.481389"
$data.Range("F8").Value = "2021-10-05 14:22:38.481391"
$data.Range("F9").Value = "2021-10-05 14:22:38.481394"
$data.Range("F10").Value = "2021-10-05 14:22:38.481397"
$data.Range("F11").Value = "2021-10-05 14:22:38.481399"
$data.Range("F12").Value = "2021-10-05 14:22:38.481402"
$data.Range("F13").Value = "2021-10-05 14:22:38.481404"
$data.Range("F14").Value = "2021-10-05 14:22:38.481407"
$data.Range("F15").Value = "2021-10-05 14:22:38.481409"
$data.Range("F16").Value = "2021-10-05 14:22:38.481412"
$data.Range("F17").Value = "2021-10-05 14:22:38.481414"
$data.Range("F18").Value = "2021-10-05 14:22:38.481417"

# Add a new "metadata" sheet right after "data" to hold panel-level info,
# refining what used to be duplicated per gene row.
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Reuse the bold/bordered/centered style from the "data" sheet header/id
# column for the header row and the A2 id cell.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$meta.Range("B2").Value = "Segmental overgrowth disorders"
$meta.Range("C2").Value = 98
$meta.Range("D2").Value = "'2.14"
$meta.Range("E2").Value = "2021-06-21T13:21:53.876767Z"
$meta.Range("F2").Value = "2021-10-05 14:22:38.477617"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/98/?format=json"

$data.Activate()
